$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the row height convention used by every other data row
$ws.Rows("27").RowHeight = 30

# Fill in the new time-sheet entry row (row 27)
$ws.Range("A27").Value = 41886
$ws.Range("B27").Value = 0.39583333333333331
$ws.Range("C27").Value = 0.47916666666666669
$ws.Range("D27").Formula = "=HOUR(C27-B27) + MINUTE(C27-B27) / 60"
$ws.Range("E27").Value = "Finishing Touches"

# Copy formatting from the previous data row (26) down into the new row (27)
$ws.Range("A26:E26").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A28").Select()
